$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
}

# Column D (Price) updates - numeric-looking strings that must stay as text
Set-TextValue $ws "D2" "241.64"
Set-TextValue $ws "D3" "21.89"
Set-TextValue $ws "D4" "5.347"
Set-TextValue $ws "D5" "0.05701"
Set-TextValue $ws "D6" "3.426"
Set-TextValue $ws "D7" "6.294"
Set-TextValue $ws "D8" "0.8056"
Set-TextValue $ws "D9" "0.8535"
Set-TextValue $ws "D11" "0.07276"
Set-TextValue $ws "D12" "0.03045"
Set-TextValue $ws "D13" "0.03145"
Set-TextValue $ws "D14" "0.09371"
Set-TextValue $ws "D15" "3.921"
Set-TextValue $ws "D16" "0.001589"
Set-TextValue $ws "D17" "0.04813"
Set-TextValue $ws "D18" "0.0005848"
Set-TextValue $ws "D19" "0.006338"

# Rows 20 and 21 swap their coin identity (Hotbit <-> BitKan)
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D20" "0.0009999"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D21" "0.004060"
$ws.Range("E21").Value = "20HotbitTokenHTB"

Set-TextValue $ws "D23" "3.718"
Set-TextValue $ws "D24" "2.170"
Set-TextValue $ws "D25" "0.3234"
Set-TextValue $ws "D27" "0.0003031"

Set-TextValue $ws "D40" "0.03822"
Set-TextValue $ws "D41" "0.006741"
Set-TextValue $ws "D42" "0.1048"
Set-TextValue $ws "D43" "0.002423"
Set-TextValue $ws "D44" "0.006514"
Set-TextValue $ws "D45" "0.00005608"
Set-TextValue $ws "D47" "0.5808"
Set-TextValue $ws "D49" "0.00002103"
